# Generate Report for Handback
#
# For the two locale sheets (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Handback DateTime" (column H) placeholder is replaced with the
#     real handback timestamp
#   - Two new columns are populated: "Latest Target File" (F) and
#     "Latest Handback File" (G) — both hyperlinked, mirroring the handoff file
#     that was produced (status is "in sync", so target/handback == the
#     original handoff files)

$wb = $excel.ActiveWorkbook

# The "Overview" sheet keeps a per-locale status summary (columns B=zh-cn,
# C=de-de) that shares the very same "Ready for handoff" string as the
# locale sheets' own Status column, so it flips to the new text too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$rows = @(
    @{
        FileName       = "59acd6e1-e3b6-4f64-b1e3-05682d531ea5.md"
        RowNum         = 2
        MdAddress      = "https://github.com/OpenLocalizationTest/oltest/blob/9ca5628870dc157c146355c0e218d1e4c6818235/e2e/59acd6e1-e3b6-4f64-b1e3-05682d531ea5.md"
    },
    @{
        FileName       = "a8c1cb60-ea24-468e-b5ce-5577592376be.md"
        RowNum         = 3
        MdAddress      = "https://github.com/OpenLocalizationTest/oltest/blob/9ca5628870dc157c146355c0e218d1e4c6818235/e2e/a8c1cb60-ea24-468e-b5ce-5577592376be.md"
    }
)

$locales = @(
    @{
        Sheet        = "zh-cn"
        HandbackTime = "2016-03-11 08:06:34"
        XlfAddress   = @{
            2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c7ec31e5a468a58949ea18ec4e1569eb22e5758/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/59acd6e1-e3b6-4f64-b1e3-05682d531ea5.a2a944104cf08ee403f468ccc1e2a94de392a8e7.zh-cn.xlf"
            3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c7ec31e5a468a58949ea18ec4e1569eb22e5758/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/a8c1cb60-ea24-468e-b5ce-5577592376be.4daf36cba7549c7aa7457a94a2a164377bfa6168.zh-cn.xlf"
        }
        XlfName      = @{
            2 = "59acd6e1-e3b6-4f64-b1e3-05682d531ea5.a2a944104cf08ee403f468ccc1e2a94de392a8e7.zh-cn.xlf"
            3 = "a8c1cb60-ea24-468e-b5ce-5577592376be.4daf36cba7549c7aa7457a94a2a164377bfa6168.zh-cn.xlf"
        }
    },
    @{
        Sheet        = "de-de"
        HandbackTime = "2016-03-11 08:06:39"
        XlfAddress   = @{
            2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7250a343020072a6592344fac94a9e48d5e45ac6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/59acd6e1-e3b6-4f64-b1e3-05682d531ea5.a2a944104cf08ee403f468ccc1e2a94de392a8e7.de-de.xlf"
            3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7250a343020072a6592344fac94a9e48d5e45ac6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/a8c1cb60-ea24-468e-b5ce-5577592376be.4daf36cba7549c7aa7457a94a2a164377bfa6168.de-de.xlf"
        }
        XlfName      = @{
            2 = "59acd6e1-e3b6-4f64-b1e3-05682d531ea5.a2a944104cf08ee403f468ccc1e2a94de392a8e7.de-de.xlf"
            3 = "a8c1cb60-ea24-468e-b5ce-5577592376be.4daf36cba7549c7aa7457a94a2a164377bfa6168.de-de.xlf"
        }
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    foreach ($row in $rows) {
        $r = $row.RowNum

        # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Range("C$r").Value = "Handed back: in sync with en-US"

        # Latest Handback DateTime
        $ws.Range("H$r").Value = $locale.HandbackTime

        # Latest Target File (F) — the handed-off markdown, now confirmed as target
        $ws.Range("F$r").Value = $row.FileName
        $ws.Hyperlinks.Add($ws.Range("F$r"), $row.MdAddress, "", "", $row.FileName) | Out-Null

        # Latest Handback File (G) — the localized xlf handed back
        $ws.Range("G$r").Value = $locale.XlfName[$r]
        $ws.Hyperlinks.Add($ws.Range("G$r"), $locale.XlfAddress[$r], "", "", $locale.XlfName[$r]) | Out-Null
    }
}
